# Albtal_R.xlsx update:
#  - the "BfFestAlbtal" label was renamed to "BahnhofsfestAlbtal"
#  - the worksheet selection moved from the whole E column (E1:E1048576)
#    down to the single cell E5, which now holds that label

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the text that used to read "BfFestAlbtal"
$ws.Range("E5").Value = "BahnhofsfestAlbtal"

# Match the new active selection recorded in the sheet view (E5 only)
$ws.Range("E5").Select() | Out-Null
